$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.10087
$ws.Range("H2").Value = 3.30261
$ws.Range("I2").Value = 0.1843884439613191
$ws.Range("J2").Value = 0.1843884439613191
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.116143
$ws.Range("N2").Value = 0.348429
$ws.Range("O2").Value = 0.01430514908838541
$ws.Range("P2").Value = 0.01430514908838541
$ws.Range("Q2").Value = 0.12785834441
$ws.Range("R2").Value = 1.15072509969
$ws.Range("S2").Value = 0.002637704181042067
$ws.Range("T2").Value = 0.002637704181042067
$ws.Range("G3").Value = 1.10087
$ws.Range("H3").Value = 3.30261
$ws.Range("I3").Value = 0.1843884439613191
$ws.Range("J3").Value = 0.1843884439613191
$ws.Range("O3").Value = 0.5605328823946109
$ws.Range("P3").Value = 0.5605328823946107
$ws.Range("Q3").Value = 5.010000656933334
$ws.Range("R3").Value = 45.0900059124
$ws.Range("S3").Value = 0.1033557859738954
$ws.Range("T3").Value = 0.1033557859738953
$ws.Range("G4").Value = 1.10087
$ws.Range("H4").Value = 3.30261
$ws.Range("I4").Value = 0.1843884439613191
$ws.Range("J4").Value = 0.1843884439613191
$ws.Range("O4").Value = 0.4251619685170038
$ws.Range("P4").Value = 0.4251619685170038
$ws.Range("Q4").Value = 3.80006563125
$ws.Range("R4").Value = 34.20059068125
$ws.Range("S4").Value = 0.07839495380638166
$ws.Range("T4").Value = 0.07839495380638166
$ws.Range("I5").Value = 0.4503925067925547
$ws.Range("J5").Value = 0.4503925067925547
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.116143
$ws.Range("N5").Value = 0.348429
$ws.Range("O5").Value = 0.01430514908838541
$ws.Range("P5").Value = 0.01430514908838541
$ws.Range("Q5").Value = 0.3123104627166667
$ws.Range("R5").Value = 2.81079416445
$ws.Range("S5").Value = 0.006442931957959132
$ws.Range("T5").Value = 0.006442931957959131
$ws.Range("I6").Value = 0.4503925067925547
$ws.Range("J6").Value = 0.4503925067925547
$ws.Range("O6").Value = 0.5605328823946109
$ws.Range("P6").Value = 0.5605328823946107
$ws.Range("S6").Value = 0.252459810041365
$ws.Range("T6").Value = 0.252459810041365
$ws.Range("I7").Value = 0.4503925067925547
$ws.Range("J7").Value = 0.4503925067925547
$ws.Range("O7").Value = 0.4251619685170038
$ws.Range("P7").Value = 0.4251619685170038
$ws.Range("S7").Value = 0.1914897647932306
$ws.Range("T7").Value = 0.1914897647932306
$ws.Range("H8").Value = 6.541494999999999
$ws.Range("I8").Value = 0.3652190492461261
$ws.Range("J8").Value = 0.3652190492461262
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.116143
$ws.Range("N8").Value = 0.348429
$ws.Range("O8").Value = 0.01430514908838541
$ws.Range("P8").Value = 0.01430514908838541
$ws.Range("Q8").Value = 0.2532496179283333
$ws.Range("R8").Value = 2.279246561355
$ws.Range("S8").Value = 0.005224512949384206
$ws.Range("T8").Value = 0.005224512949384207
$ws.Range("H9").Value = 6.541494999999999
$ws.Range("I9").Value = 0.3652190492461261
$ws.Range("J9").Value = 0.3652190492461262
$ws.Range("O9").Value = 0.5605328823946109
$ws.Range("P9").Value = 0.5605328823946107
$ws.Range("Q9").Value = 9.923331621755555
$ws.Range("R9").Value = 89.3099845958
$ws.Range("S9").Value = 0.2047172863793504
$ws.Range("T9").Value = 0.2047172863793504
$ws.Range("H10").Value = 6.541494999999999
$ws.Range("I10").Value = 0.3652190492461261
$ws.Range("J10").Value = 0.3652190492461262
$ws.Range("O10").Value = 0.4251619685170038
$ws.Range("P10").Value = 0.4251619685170038
$ws.Range("Q10").Value = 7.526807684374998
$ws.Range("R10").Value = 67.74126915937499
$ws.Range("S10").Value = 0.1552772499173916
$ws.Range("T10").Value = 0.1552772499173916
